# Generate Report for Handback
# Updates the localization-status workbook to reflect that the zh-cn and
# de-de handoffs have been handed back and are back in sync with en-US.

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"
$hyperlinkColor = 15570276   # matches the workbook's existing HyperLink font color (FF6495ED)

# --- Overview sheet: status cells mirror the per-language sheets ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2").Value = $statusText
$wsOverview.Range("C2").Value = $statusText
$wsOverview.Range("B3").Value = $statusText
$wsOverview.Range("C3").Value = $statusText

# --- Per-language sheets: zh-cn and de-de ---
# Row 2 -> 27191251-2634-40f4-a04b-7b7a3a1b317e.* (source .md + its zh-cn/de-de .xlf)
# Row 3 -> eb5b6d7f-1b9b-4f04-88e3-fae0653135b2.* (source .md + its zh-cn/de-de .xlf)

$rows = @(
    @{
        Row = 2
        SourceName   = "27191251-2634-40f4-a04b-7b7a3a1b317e.md"
        SourceUrl    = "https://github.com/OpenLocalizationTest/oltest/blob/1957997fc0b2041ce8f6924bad1062ba4d886db6/e2e/27191251-2634-40f4-a04b-7b7a3a1b317e.md"
        HandoffNameZh = "27191251-2634-40f4-a04b-7b7a3a1b317e.4e9d883f10139ba1bc316ffd214f987f2fd835fa.zh-cn.xlf"
        HandoffUrlZh  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/11fc754cf784e6ae89328d77cdbd2ff018f0eff0/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/27191251-2634-40f4-a04b-7b7a3a1b317e.4e9d883f10139ba1bc316ffd214f987f2fd835fa.zh-cn.xlf"
        HandoffNameDe = "27191251-2634-40f4-a04b-7b7a3a1b317e.4e9d883f10139ba1bc316ffd214f987f2fd835fa.de-de.xlf"
        HandoffUrlDe  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0edbf6406e77911462ac2b6839c6b9e5f0328bd1/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/27191251-2634-40f4-a04b-7b7a3a1b317e.4e9d883f10139ba1bc316ffd214f987f2fd835fa.de-de.xlf"
    },
    @{
        Row = 3
        SourceName   = "eb5b6d7f-1b9b-4f04-88e3-fae0653135b2.md"
        SourceUrl    = "https://github.com/OpenLocalizationTest/oltest/blob/1957997fc0b2041ce8f6924bad1062ba4d886db6/e2e/eb5b6d7f-1b9b-4f04-88e3-fae0653135b2.md"
        HandoffNameZh = "eb5b6d7f-1b9b-4f04-88e3-fae0653135b2.c5f9a66ec1bfd70a8e49530b4826e8162b806280.zh-cn.xlf"
        HandoffUrlZh  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/11fc754cf784e6ae89328d77cdbd2ff018f0eff0/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/eb5b6d7f-1b9b-4f04-88e3-fae0653135b2.c5f9a66ec1bfd70a8e49530b4826e8162b806280.zh-cn.xlf"
        HandoffNameDe = "eb5b6d7f-1b9b-4f04-88e3-fae0653135b2.c5f9a66ec1bfd70a8e49530b4826e8162b806280.de-de.xlf"
        HandoffUrlDe  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0edbf6406e77911462ac2b6839c6b9e5f0328bd1/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/eb5b6d7f-1b9b-4f04-88e3-fae0653135b2.c5f9a66ec1bfd70a8e49530b4826e8162b806280.de-de.xlf"
    }
)

$langs = @(
    @{ Sheet = "zh-cn"; HandbackTime = "2016-03-09 03:41:31"; NameKey = "HandoffNameZh"; UrlKey = "HandoffUrlZh" },
    @{ Sheet = "de-de"; HandbackTime = "2016-03-09 03:42:09"; NameKey = "HandoffNameDe"; UrlKey = "HandoffUrlDe" }
)

foreach ($lang in $langs) {
    $ws = $wb.Worksheets.Item($lang.Sheet)
    $handbackTime = $lang.HandbackTime

    foreach ($r in $rows) {
        $row = $r.Row
        $eCell = "E$row"
        $fCell = "F$row"
        $gCell = "G$row"

        $handoffName = $r[$lang.NameKey]
        $handoffUrl = $r[$lang.UrlKey]

        # Status column
        $ws.Range("B$row").Value = $statusText

        # Latest Target File (E) - same source file that was handed back
        $ws.Range($eCell).Value = $r.SourceName
        $ws.Hyperlinks.Add($ws.Range($eCell), $r.SourceUrl, "", "", $r.SourceName)
        $ws.Range($eCell).Font.Underline = 2
        $ws.Range($eCell).Font.Color = $hyperlinkColor

        # Latest Handback File (F) - same handoff xlf file, now handed back
        $ws.Range($fCell).Value = $handoffName
        $ws.Hyperlinks.Add($ws.Range($fCell), $handoffUrl, "", "", $handoffName)
        $ws.Range($fCell).Font.Underline = 2
        $ws.Range($fCell).Font.Color = $hyperlinkColor

        # Latest Handback DateTime (G)
        $ws.Range($gCell).Value = $handbackTime
    }
}
